# Updates the cryptocurrency symbol/price table on the active sheet to
# reflect the latest scrape (GitHub Actions run on Sun Feb 12 07:08:55 UTC 2023).
#
# Each entry is [CellRef, NewValue, LooksNumeric]. Cells whose new text looks
# like a plain number or a percentage ("307.39", "-0.40%", "7", ...) must be
# written with the cell's NumberFormat forced to Text ("@") first; otherwise
# Excel auto-converts the literal text into a numeric/percentage value instead
# of keeping it as the original string the source data used.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    ,@("D2", "307.39", $true)
    ,@("E2", "-0.40%", $true)
    ,@("G2", "7", $true)
    ,@("D3", "41.06", $true)
    ,@("E3", "0.69%", $true)
    ,@("G3", "7", $true)
    ,@("D4", "5.246", $true)
    ,@("E4", "2.54%", $true)
    ,@("G4", "7", $true)
    ,@("D5", "0.07665", $true)
    ,@("E5", "0.65%", $true)
    ,@("G5", "7", $true)
    ,@("D6", "1.636", $true)
    ,@("E6", "0.78%", $true)
    ,@("G6", "7", $true)
    ,@("D7", "0.9151", $true)
    ,@("E7", "1.54%", $true)
    ,@("G7", "7", $true)
    ,@("D8", "2.443", $true)
    ,@("E8", "0.57%", $true)
    ,@("G8", "7", $true)
    ,@("E9", "14.27%", $true)
    ,@("G9", "7", $true)
    ,@("D10", "0.1825", $true)
    ,@("E10", "3.57%", $true)
    ,@("G10", "7", $true)
    ,@("D11", "0.09151", $true)
    ,@("E11", "-0.02%", $true)
    ,@("G11", "7", $true)
    ,@("D12", "0.04268", $true)
    ,@("E12", "2.15%", $true)
    ,@("G12", "7", $true)
    ,@("E13", "0.01%", $true)
    ,@("G13", "7", $true)
    ,@("E14", "0.54%", $true)
    ,@("G14", "7", $true)
    ,@("D15", "0.005820", $true)
    ,@("E15", "-0.54%", $true)
    ,@("G15", "7", $true)
    ,@("G16", "7", $true)
    ,@("B17", "HotbitToken", $false)
    ,@("C17", "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb", $false)
    ,@("D17", "0.004278", $true)
    ,@("E17", "4.61%", $true)
    ,@("G17", "7", $true)
    ,@("B18", "LEO", $false)
    ,@("C18", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo", $false)
    ,@("D18", "3.346", $true)
    ,@("E18", "-0.21%", $true)
    ,@("G18", "7", $true)
    ,@("B19", "GateToken", $false)
    ,@("C19", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt", $false)
    ,@("D19", "4.308", $true)
    ,@("E19", "1.37%", $true)
    ,@("G19", "7", $true)
    ,@("B20", "BitpandaEcosystemToken", $false)
    ,@("C20", "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best", $false)
    ,@("D20", "0.3336", $true)
    ,@("E20", "1.13%", $true)
    ,@("G20", "7", $true)
    ,@("B21", "MCDex", $false)
    ,@("C21", "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb", $false)
    ,@("D21", "7.317", $true)
    ,@("E21", "11.83%", $true)
    ,@("G21", "7", $true)
    ,@("B22", "ProBitToken", $false)
    ,@("C22", "https://coinranking.com/coin/lQP4d6T2+probittoken-prob", $false)
    ,@("D22", "0.1384", $true)
    ,@("E22", "1.40%", $true)
    ,@("G22", "7", $true)
    ,@("B23", "ZBToken", $false)
    ,@("C23", "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb", $false)
    ,@("D23", "0.2715", $true)
    ,@("E23", "1.27%", $true)
    ,@("G23", "7", $true)
    ,@("B24", "CoinExToken", $false)
    ,@("C24", "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet", $false)
    ,@("D24", "0.04063", $true)
    ,@("E24", "-0.29%", $true)
    ,@("G24", "7", $true)
    ,@("B25", "BitKan", $false)
    ,@("C25", "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan", $false)
    ,@("D25", "0.001264", $true)
    ,@("E25", "3.42%", $true)
    ,@("G25", "7", $true)
    ,@("E26", "-2.14%", $true)
    ,@("G26", "7", $true)
    ,@("G27", "7", $true)
    ,@("G28", "7", $true)
    ,@("G29", "7", $true)
    ,@("G30", "7", $true)
    ,@("G31", "7", $true)
    ,@("G32", "7", $true)
    ,@("G33", "7", $true)
    ,@("G34", "7", $true)
    ,@("G35", "7", $true)
    ,@("G36", "7", $true)
    ,@("G37", "7", $true)
    ,@("D38", "0.02478", $true)
    ,@("E38", "4.05%", $true)
    ,@("G38", "7", $true)
    ,@("D39", "0.05314", $true)
    ,@("E39", "2.56%", $true)
    ,@("G39", "7", $true)
    ,@("D40", "0.007847", $true)
    ,@("E40", "1.06%", $true)
    ,@("G40", "7", $true)
    ,@("D41", "0.1312", $true)
    ,@("E41", "0.89%", $true)
    ,@("G41", "7", $true)
    ,@("D42", "0.006880", $true)
    ,@("E42", "1.69%", $true)
    ,@("G42", "7", $true)
    ,@("E43", "-1.89%", $true)
    ,@("G43", "7", $true)
    ,@("D44", "0.007628", $true)
    ,@("E44", "-8.45%", $true)
    ,@("G44", "7", $true)
    ,@("D45", "0.3065", $true)
    ,@("E45", "-0.45%", $true)
    ,@("G45", "7", $true)
    ,@("D46", "0.00006724", $true)
    ,@("E46", "-3.26%", $true)
    ,@("G46", "7", $true)
    ,@("E47", "0.17%", $true)
    ,@("G47", "7", $true)
    ,@("D48", "0.1699", $true)
    ,@("E48", "430.93%", $true)
    ,@("G48", "7", $true)
    ,@("D49", "0.003106", $true)
    ,@("E49", "-26.08%", $true)
    ,@("G49", "7", $true)
    ,@("E50", "0.17%", $true)
    ,@("G50", "7", $true)
    ,@("E51", "0.17%", $true)
    ,@("G51", "7", $true)
)

foreach ($u in $updates) {
    $cellRef = $u[0]
    $newValue = $u[1]
    $looksNumeric = $u[2]

    $rng = $ws.Range($cellRef)
    if ($looksNumeric) {
        $rng.NumberFormat = "@"
    }
    $rng.Value = $newValue
}

Write-Output "Applied $($updates.Count) cell updates"
